$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.083535801138282
$ws.Range("C2").Value = 0.2991360873785709
$ws.Range("D2").Value = 0.02766834803760077
$ws.Range("F2").Value = 0.9256811883042673
$ws.Range("G2").Value = 0.7780272031377962
$ws.Range("H2").Value = 0.8312015576536709
$ws.Range("I2").Value = 0.7994050612249666
$ws.Range("L2").Value = 0.27140343584108
$ws.Range("M2").Value = 0.2517476747808516
$ws.Range("B3").Value = 0.9786524220576212
$ws.Range("C3").Value = 0.2757984692949265
$ws.Range("D3").Value = 0.02728259605396843
$ws.Range("F3").Value = 0.9180137608706929
$ws.Range("G3").Value = 0.7713470043885309
$ws.Range("H3").Value = 0.8339881759197567
$ws.Range("I3").Value = 0.8070507183222588
$ws.Range("L3").Value = 0.2700854042169496
$ws.Range("M3").Value = 0.2361886505213775
$ws.Range("B4").Value = 0.9143940898700862
$ws.Range("C4").Value = 0.2613585856882992
$ws.Range("D4").Value = 0.02704863232600374
$ws.Range("F4").Value = 0.9140676193005959
$ws.Range("G4").Value = 0.7679697337431008
$ws.Range("H4").Value = 0.8362446933044367
$ws.Range("I4").Value = 0.8123402092025742
$ws.Range("L4").Value = 0.2694203457237094
$ws.Range("M4").Value = 0.22672315716882
$ws.Range("B5").Value = 0.8882446368543242
$ws.Range("C5").Value = 0.2554465471264109
$ws.Range("D5").Value = 0.02695402986009654
$ws.Range("F5").Value = 0.9126504674355758
$ws.Range("G5").Value = 0.7667748563823267
$ws.Range("H5").Value = 0.8373011880083538
$ws.Range("I5").Value = 0.8146449868697125
$ws.Range("L5").Value = 0.2691856033814801
$ws.Range("M5").Value = 0.2228881035842889
$ws.Range("B6").Value = 0.8839047651313194
$ws.Range("C6").Value = 0.2544631899856995
$ws.Range("D6").Value = 0.02693836625592638
$ws.Range("F6").Value = 0.9124266647596428
$ws.Range("G6").Value = 0.7665873807962953
$ws.Range("H6").Value = 0.8374848825872476
$ws.Range("I6").Value = 0.8150366983545716
$ws.Range("L6").Value = 0.2691488157739457
$ws.Range("M6").Value = 0.2222526413234505
$ws.Range("B7").Value = 0.9140412807092844
$ws.Range("C7").Value = 0.2612789657387964
$ws.Range("D7").Value = 0.02704735347248644
$ws.Range("F7").Value = 0.9140477348147797
$ws.Range("G7").Value = 0.7679528857981666
$ws.Range("H7").Value = 0.836258387384234
$ws.Range("I7").Value = 0.8123706883411437
$ws.Range("L7").Value = 0.2694170330193089
$ws.Range("M7").Value = 0.22667134618424
$ws.Range("B8").Value = 1.047343423019413
$ws.Range("C8").Value = 0.2911122708147218
$ws.Range("D8").Value = 0.0275347486047508
$ws.Range("F8").Value = 0.9228790043205919
$ws.Range("G8").Value = 0.7755730553871132
$ws.Range("H8").Value = 0.8320490199914872
$ws.Range("I8").Value = 0.8019175710568689
$ws.Range("L8").Value = 0.2709190504914645
$ws.Range("M8").Value = 0.2463647708374879
$ws.Range("B9").Value = 1.309832892266229
$ws.Range("C9").Value = 0.3487363214255481
$ws.Range("D9").Value = 0.02851295455731062
$ws.Range("F9").Value = 0.9462701819887229
$ws.Range("G9").Value = 0.7963005787861306
$ws.Range("H9").Value = 0.8281341804122349
$ws.Range("I9").Value = 0.7861566395756441
$ws.Range("L9").Value = 0.2750089300318379
$ws.Range("M9").Value = 0.2856767554272679
$ws.Range("B10").Value = 1.503325126266816
$ws.Range("C10").Value = 0.3905379559373898
$ws.Range("D10").Value = 0.02924474478306038
$ws.Range("F10").Value = 0.9672020984542939
$ws.Range("G10").Value = 0.815108948342413
$ws.Range("H10").Value = 0.8279203378458959
$ws.Range("I10").Value = 0.7774881875287676
$ws.Range("L10").Value = 0.2787124026937846
$ws.Range("M10").Value = 0.3149805340152767
$ws.Range("B11").Value = 1.591486373694579
$ws.Range("C11").Value = 0.4094390246575585
$ws.Range("D11").Value = 0.02958039232466092
$ws.Range("F11").Value = 0.9775477623253579
$ws.Range("G11").Value = 0.8244543765600696
$ws.Range("H11").Value = 0.8284048749521276
$ws.Range("I11").Value = 0.7741817708115661
$ws.Range("L11").Value = 0.2805491775460212
$ws.Range("M11").Value = 0.3284030203531927
$ws.Range("B12").Value = 1.624890322857368
$ws.Range("C12").Value = 0.4165798299752055
$ws.Range("D12").Value = 0.02970787893985261
$ws.Range("F12").Value = 0.9815845481800096
$ws.Range("G12").Value = 0.8281076481745799
$ws.Range("H12").Value = 0.8286723028431311
$ws.Range("I12").Value = 0.7730217057961326
$ws.Range("L12").Value = 0.2812665877313236
$ws.Range("M12").Value = 0.3334989317078296
$ws.Range("B13").Value = 1.617695347062352
$ws.Range("C13").Value = 0.4150426708778525
$ws.Range("D13").Value = 0.02968040551041184
$ws.Range("F13").Value = 0.9807098461041051
$ws.Range("G13").Value = 0.8273157510857629
$ws.Range("H13").Value = 0.8286109695544894
$ws.Range("I13").Value = 0.7732674480412669
$ws.Range("L13").Value = 0.2811111082373401
$ws.Range("M13").Value = 0.332400855114777
$ws.Range("B14").Value = 1.594234157760241
$ws.Range("C14").Value = 0.4100268367847661
$ws.Range("D14").Value = 0.02959087309102415
$ws.Range("F14").Value = 0.9778774801956018
$ws.Range("G14").Value = 0.8247526366146758
$ws.Range("H14").Value = 0.8284251924328316
$ws.Range("I14").Value = 0.7740844852936704
$ws.Range("L14").Value = 0.2806077611410132
$ws.Range("M14").Value = 0.3288220020949097
$ws.Range("B15").Value = 1.579865971847653
$ws.Range("C15").Value = 0.4069523256193577
$ws.Range("D15").Value = 0.02953608159054255
$ws.Range("F15").Value = 0.9761581055015256
$ws.Range("G15").Value = 0.8231975732174845
$ws.Range("H15").Value = 0.8283223389146883
$ws.Range("I15").Value = 0.7745969381185489
$ws.Range("L15").Value = 0.2803022937586093
$ws.Range("M15").Value = 0.3266315544101488
$ws.Range("B16").Value = 1.49756632347777
$ws.Range("C16").Value = 0.3893004077456794
$ws.Range("D16").Value = 0.02922286381108208
$ws.Range("F16").Value = 0.9665426158470751
$ws.Range("G16").Value = 0.814514159752548
$ws.Range("H16").Value = 0.8279004037233335
$ws.Range("I16").Value = 0.7777171214461447
$ws.Range("L16").Value = 0.2785954246339912
$ws.Range("M16").Value = 0.314105184111817
$ws.Range("B17").Value = 1.447113443240482
$ws.Range("C17").Value = 0.3784420593496236
$ws.Range("D17").Value = 0.02903141142709487
$ws.Range("F17").Value = 0.9608552528917045
$ws.Range("G17").Value = 0.8093899541829472
$ws.Range("H17").Value = 0.8277907758926233
$ws.Range("I17").Value = 0.7797946762947134
$ws.Range("L17").Value = 0.2775872578550889
$ws.Range("M17").Value = 0.3064441448860293
$ws.Range("B18").Value = 1.418107580782191
$ws.Range("C18").Value = 0.3721858352634797
$ws.Range("D18").Value = 0.0289215528249791
$ws.Range("F18").Value = 0.957661516021048
$ws.Range("G18").Value = 0.8065168988478035
$ws.Range("H18").Value = 0.8277824672665162
$ws.Range("I18").Value = 0.7810495545099414
$ws.Range("L18").Value = 0.2770216983491736
$ws.Range("M18").Value = 0.3020463895587895
$ws.Range("B19").Value = 1.408289010056251
$ws.Range("C19").Value = 0.3700657360736841
$ws.Range("D19").Value = 0.02888440156326411
$ws.Range("F19").Value = 0.9565934590983574
$ws.Range("G19").Value = 0.8055568588110731
$ws.Range("H19").Value = 0.8277890483463182
$ws.Range("I19").Value = 0.7814847139543986
$ws.Range("L19").Value = 0.2768326676356878
$ws.Range("M19").Value = 0.3005588785550231
$ws.Range("B20").Value = 1.452482869504536
$ws.Range("C20").Value = 0.3795990666421005
$ws.Range("D20").Value = 0.02905176506143192
$ws.Range("F20").Value = 0.9614526582366238
$ws.Range("G20").Value = 0.8099277436672168
$ws.Range("H20").Value = 0.8277967778849131
$ws.Range("I20").Value = 0.7795673121405571
$ws.Range("L20").Value = 0.2776930977643133
$ws.Range("M20").Value = 0.3072587790388042
$ws.Range("B21").Value = 1.601124764837039
$ws.Range("C21").Value = 0.411500560861441
$ws.Range("D21").Value = 0.02961716059730435
$ws.Range("F21").Value = 0.9787061768340948
$ws.Range("G21").Value = 0.8255023750132722
$ws.Range("H21").Value = 0.8284774792201404
$ws.Range("I21").Value = 0.7738420013812828
$ws.Range("L21").Value = 0.2807550130214054
$ws.Range("M21").Value = 0.3298728432770233
$ws.Range("B22").Value = 1.698382273713378
$ws.Range("C22").Value = 0.4322531066660531
$ws.Range("D22").Value = 0.02998891342439691
$ws.Range("F22").Value = 0.9906769093553009
$ws.Range("G22").Value = 0.8363482711009169
$ws.Range("H22").Value = 0.8294117784276978
$ws.Range("I22").Value = 0.7706366388093286
$ws.Range("L22").Value = 0.2828835883161105
$ws.Range("M22").Value = 0.3447287813912396
$ws.Range("B23").Value = 1.646464249157304
$ws.Range("C23").Value = 0.4211859913557703
$ws.Range("D23").Value = 0.0297903012570444
$ws.Range("F23").Value = 0.9842241302282702
$ws.Range("G23").Value = 0.8304983048862198
$ws.Range("H23").Value = 0.8288682505766189
$ws.Range("I23").Value = 0.772298177269505
$ws.Range("L23").Value = 0.2817358672782291
$ws.Range("M23").Value = 0.3367929446699947
$ws.Range("B24").Value = 1.450055351287631
$ws.Range("C24").Value = 0.3790760260511945
$ws.Range("D24").Value = 0.02904256252826443
$ws.Range("F24").Value = 0.9611823345787798
$ws.Range("G24").Value = 0.8096843819651838
$ws.Range("H24").Value = 0.8277938939543219
$ws.Range("I24").Value = 0.7796699151563899
$ws.Range("L24").Value = 0.2776452037812192
$ws.Range("M24").Value = 0.3068904621278676
$ws.Range("B25").Value = 1.238708571902293
$ws.Range("C25").Value = 0.333241489862786
$ws.Range("D25").Value = 0.02824598346405693
$ws.Range("F25").Value = 0.9392875237934106
$ws.Range("G25").Value = 0.7900684399191107
$ws.Range("H25").Value = 0.8287269195203208
$ws.Range("I25").Value = 0.7899108261848795
$ws.Range("L25").Value = 0.2737798922775809
$ws.Range("M25").Value = 0.2749677563406365
